# Add new columns C ("kpi PPR 1/3"), D ("anslag kpi PRR 2/22") and
# E ("Inflasjonsmål") of data next to the existing A (dato) / B (styringsrent)
# columns on Ark1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -------------------------------------------------------
$ws.Range("C1").Value = "kpi PPR 1/3"
$ws.Range("D1").Value = "anslag kpi PRR 2/22"
$ws.Range("E1").Value = "Inflasjonsmål"

# ---- Column widths ------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 16
$ws.Columns.Item(5).ColumnWidth = 11

# ---- Data rows ----------------------------------------------------------
# columns: row, C (kpi), D (anslag, $null if blank), E (inflasjonsmaal)
$data = @(
    @(2,  1.1499999999999999, $null, 2),
    @(3,  1.1399999999999999, $null, 2),
    @(4,  1.56,               $null, 2),
    @(5,  1.29,               $null, 2),
    @(6,  2.97,               $null, 2),
    @(7,  2.83,               $null, 2),
    @(8,  3.49,               $null, 2),
    @(9,  4.6399999999999997, $null, 2),
    @(10, 3.81,               $null, 2),
    @(11, 5.82,               $null, 2),
    @(12, 6.74,               $null, 2),
    @(13, 6.63,               6.66,  2),
    @(14, 6.44,               6.64,  2),
    @(15, 5.14,               5.49,  2),
    @(16, 4.24,               4.04,  2),
    @(17, 3.78,               3.15,  2),
    @(18, 3.95,               3.12,  2),
    @(19, 3.67,               2.95,  2),
    @(20, 3.03,               2.63,  2),
    @(21, 2.71,               2.63,  2)
)

foreach ($item in $data) {
    $row = $item[0]
    $kpi = $item[1]
    $anslag = $item[2]
    $inflation = $item[3]

    # Column C picks up the same number format / font colour as column B
    # on this row (black through row 13, blue "forecast" colour from row 14).
    $ws.Range("B$row").Copy()
    $ws.Range("C$row").PasteSpecial(-4122)
    $ws.Range("C$row").Value = $kpi

    # Column D: rows 2-12 are blank but still carry the black number style
    # (copied from B2); rows 13-21 hold a value styled like the blue
    # "forecast" cells (copied from B14).
    if ($row -le 12) {
        $ws.Range("B2").Copy()
        $ws.Range("D$row").PasteSpecial(-4122)
    } else {
        $ws.Range("B14").Copy()
        $ws.Range("D$row").PasteSpecial(-4122)
        $ws.Range("D$row").Value = $anslag
    }

    # Column E always uses the plain black number style.
    $ws.Range("B2").Copy()
    $ws.Range("E$row").PasteSpecial(-4122)
    $ws.Range("E$row").Value = $inflation
}

$excel.CutCopyMode = 0
$ws.Range("F7").Select()
